# Update "想去人数" (attendance interest count) values in the 展览 and 全部类型 sheets.
$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 265
$wsExhibit.Range("F3").Value = 85
$wsExhibit.Range("F4").Value = 949
$wsExhibit.Range("F5").Value = 545

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 265
$wsAll.Range("F3").Value = 85
$wsAll.Range("F4").Value = 949
$wsAll.Range("F6").Value = 545
